$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of new shared strings must be: "Foreign Key", "USER_ID", "Foreign key"
# Add new fields to the "Users" table (A/B/C columns) - row 15
$ws.Range("C15").Value = "Foreign Key"
$ws.Range("A15").Value = "USER_ID"

# Add new fields to the "Pucharses" table (F/G/H columns) - row 8
$ws.Range("F8").Value = "USER_ID"
$ws.Range("H8").Value = "Foreign key"

# Match the style of A15 to the style used by A13/A14 (the "MUL"/Key marker style)
$ws.Range("A15").Style = $ws.Range("A13").Style

# Update the active cell selection
$ws.Range("J9").Select()
